# BGUSO-93 WPFGen: Apply conditional visibility : In progress
# Adds a new "NotifyPropertyChanged" worksheet that derives INotifyPropertyChanged
# boilerplate (backing field + accessor) from a "Modifier Type PropName { get; set; }"
# style property declaration typed into column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new worksheet as the last tab (after WPF_temp_test)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "NotifyPropertyChanged"

# ---------------------------------------------------------------------------
# 2) Rows 3-7: the five original property declarations (entered as a block,
#    hence the shared formulas across B3:B7, D3:D7, etc.)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "        public string NameUkr { get; set; }"
$ws.Range("A4").Value = "        public CountryInfo JurisdictionCountry { get; set; }"
$ws.Range("A5").Value = "        public string CourtRegion { get; set; }"
$ws.Range("A6").Value = "        public string CourtID { get; set; }"
$ws.Range("A7").Value = "        public CourtInstanceType Instance { get; set; }"

$ws.Range("B3:B7").Formula = "=TRIM(A3)"

$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1

$ws.Range("D3:D7").Formula = '=FIND(" ",$B3)'
$ws.Range("E3:F7").Formula = '=FIND(" ",$B3,D3+1)'

$ws.Range("G3:G7").Formula = '=TRIM(MID($B3,C3,D3-C3))'
$ws.Range("H3:H7").Formula = '=TRIM(MID($B3,D3,E3-D3))'
$ws.Range("I3:I7").Formula = '=TRIM(MID($B3,E3,F3-E3))'
$ws.Range("J3:J7").Formula = '="_"&I3'
$ws.Range("K3:K7").Formula = '="private " & H3 & " " & J3 & ";"'
$ws.Range("L3:L7").Formula = '=G3& " " &H3& " " &I3 & " { get { return " & J3 & "; } set { " &J3 & " = value; OnPropertyChanged(" & CHAR(34) & I3 & CHAR(34) & "); } }"'

# ---------------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Type"
$ws.Range("G1").Value = "Modifier"
$ws.Range("I1").Value = "PropName"
$ws.Range("J1").Value = "FieldName"
$ws.Range("K1").Value = "FieldDecl"
$ws.Range("L1").Value = "Accessor"
$ws.Range("G1:M1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4) Row 2: ShortTermRatingValueOther, added on its own afterwards (so it is
#    NOT part of the shared-formula groups used by rows 3-7)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "public string ShortTermRatingValueOther { get; set; }"
$ws.Range("B2").Formula = "=TRIM(A2)"
$ws.Range("C2").Value = 1
$ws.Range("D2").Formula = '=FIND(" ",$B2)'
$ws.Range("E2").Formula = '=FIND(" ",$B2,D2+1)'
$ws.Range("F2").Formula = '=FIND(" ",$B2,E2+1)'
$ws.Range("G2").Formula = '=TRIM(MID($B2,C2,D2-C2))'
$ws.Range("H2").Formula = '=TRIM(MID($B2,D2,E2-D2))'
$ws.Range("I2").Formula = '=TRIM(MID($B2,E2,F2-E2))'
$ws.Range("J2").Formula = '="_"&I2'
$ws.Range("K2").Formula = '="private " & H2 & " " & J2 & ";"'
$ws.Range("L2").Formula = '=G2& " " &H2& " " &I2 & " { get { return " & J2 & "; } set { " &J2 & " = value; OnPropertyChanged(" & CHAR(34) & I2 & CHAR(34) & "); } }"'

# ---------------------------------------------------------------------------
# 5) Helper columns B:J are hidden - only A (source), K (field) and L
#    (accessor) are meant to be read
# ---------------------------------------------------------------------------
$helperCols = $ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item(1, 10)).EntireColumn
$helperCols.ColumnWidth = 0
$helperCols.Hidden = $true

# ---------------------------------------------------------------------------
# 6) View state: select L2, make this the active/visible tab
# ---------------------------------------------------------------------------
[void]$ws.Range("L2").Select()
$ws.Activate()
